$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1108.7567
$ws.Range("I28").Value = 345.96155
$ws.Range("K28").Value = 345.96155
$ws.Range("M28").Value = 139.03845
$ws.Range("H98").Value = 1566.7778
$ws.Range("I98").Value = 1566.7778
$ws.Range("K98").Value = 1566.7778
$ws.Range("M98").Value = -68.77780000000007
$ws.Range("H107").Value = 23881300
$ws.Range("J107").Value = 167077.67
$ws.Range("L107").Value = 167077.67
$ws.Range("N107").Value = -170917.67
$ws.Range("H122").Value = 1566.7778
$ws.Range("I122").Value = 1566.7778
$ws.Range("K122").Value = 4700.3334
$ws.Range("M122").Value = -2250.3334
$ws.Range("H129").Value = 26318064
$ws.Range("I129").Value = 45455956
$ws.Range("K129").Value = 136367868
$ws.Range("M129").Value = -136362868
$ws.Range("H135").Value = 1063.4857
$ws.Range("I135").Value = 651.3182
$ws.Range("J135").Value = 1761
$ws.Range("K135").Value = 5861.8638
$ws.Range("L135").Value = 15849
$ws.Range("M135").Value = -3326.8638
$ws.Range("N135").Value = -20919
$ws.Range("H138").Value = 2883.1448
$ws.Range("J138").Value = 3562.125
$ws.Range("L138").Value = 10686.375
$ws.Range("N138").Value = -20966.375
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = 0
$ws.Range("H92").Value = 41474.5
$ws.Range("J92").Value = 41474.5
$ws.Range("L92").Value = 41474.5
$ws.Range("N92").Value = -46466.5
$ws.Range("H102").Value = 3475260.5
$ws.Range("I102").Value = 4389198.5
$ws.Range("K102").Value = 4389198.5
$ws.Range("M102").Value = -4387576.5
$ws.Range("H122").Value = 2452261.2
$ws.Range("I122").Value = 3097938
$ws.Range("K122").Value = 9293814
$ws.Range("M122").Value = -9291364
$ws.Range("H132").Value = 25615.512
$ws.Range("I132").Value = 5954.227
$ws.Range("K132").Value = 17862.681
$ws.Range("M132").Value = -15332.681
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 25499476
$ws.Range("I86").Value = 54169600
$ws.Range("J86").Value = 14919
$ws.Range("K86").Value = 54169600
$ws.Range("L86").Value = 14919
$ws.Range("M86").Value = -54168477
$ws.Range("N86").Value = -17165
$ws.Range("H89").Value = 25499476
$ws.Range("I89").Value = 54169600
$ws.Range("J89").Value = 14919
$ws.Range("K89").Value = 270848000
$ws.Range("L89").Value = 74595
$ws.Range("M89").Value = -270842384
$ws.Range("N89").Value = -85827
$ws.Range("H105").Value = 10421008
$ws.Range("I105").Value = 12504610
$ws.Range("K105").Value = 12504610
$ws.Range("M105").Value = -12502863
$ws.Range("H132").Value = 89000
$ws.Range("J132").Value = 89000
$ws.Range("L132").Value = 89000
$ws.Range("N132").Value = -99120
$ws.Range("H134").Value = 14557.125
$ws.Range("I134").Value = 12083.368
$ws.Range("J134").Value = 23957.4
$ws.Range("K134").Value = 36250.104
$ws.Range("L134").Value = 71872.20000000001
$ws.Range("M134").Value = -33715.104
$ws.Range("N134").Value = -76942.20000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6057.6665
$ws.Range("I58").Value = 7075.0586
$ws.Range("J58").Value = 3586.8572
$ws.Range("K58").Value = 7075.0586
$ws.Range("L58").Value = 3586.8572
$ws.Range("M58").Value = -6872.0586
$ws.Range("N58").Value = -3992.8572
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26996
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84984
$ws.Range("H94").Value = 2569.2
$ws.Range("I94").Value = 3333
$ws.Range("K94").Value = 3333
$ws.Range("M94").Value = -2882
$ws.Range("H134").Value = 8400.5625
$ws.Range("I134").Value = 5805.609
$ws.Range("K134").Value = 17416.827
$ws.Range("M134").Value = -14881.827
$ws.Range("H136").Value = 6057.6665
$ws.Range("I136").Value = 7075.0586
$ws.Range("J136").Value = 3586.8572
$ws.Range("K136").Value = 21225.1758
$ws.Range("L136").Value = 10760.5716
$ws.Range("M136").Value = -18675.1758
$ws.Range("N136").Value = -15860.5716
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 50211.11
$ws.Range("J37").Value = 50211.11
$ws.Range("L37").Value = 150633.33
$ws.Range("N37").Value = -150857.33
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5908.241
$ws.Range("I132").Value = 4846.923
$ws.Range("J132").Value = 8086.737
$ws.Range("K132").Value = 14540.769
$ws.Range("L132").Value = 24260.211
$ws.Range("M132").Value = -12010.769
$ws.Range("N132").Value = -29320.211
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 25416.666
$ws.Range("J51").Value = 25416.666
$ws.Range("L51").Value = 25416.666
$ws.Range("N51").Value = -26372.666
$ws.Range("H53").Value = 21680.5
$ws.Range("J53").Value = 21680.5
$ws.Range("L53").Value = 21680.5
$ws.Range("N53").Value = -22716.5
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = 0
$ws.Range("H82").Value = 2138851.8
$ws.Range("I82").Value = 3474647
$ws.Range("J82").Value = 1579.4
$ws.Range("K82").Value = 3474647
$ws.Range("L82").Value = 1579.4
$ws.Range("M82").Value = -3474286
$ws.Range("N82").Value = -2301.4
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = 0
$ws.Range("H85").Value = 2138851.8
$ws.Range("I85").Value = 3474647
$ws.Range("J85").Value = 1579.4
$ws.Range("K85").Value = 3474647
$ws.Range("L85").Value = 1579.4
$ws.Range("M85").Value = -3473399
$ws.Range("N85").Value = -4075.4
$ws.Range("H88").Value = 24518
$ws.Range("I88").Value = 24518
$ws.Range("K88").Value = 24518
$ws.Range("M88").Value = -24090
$ws.Range("H91").Value = 24518
$ws.Range("I91").Value = 24518
$ws.Range("K91").Value = 24518
$ws.Range("M91").Value = -23036
$ws.Range("H106").Value = 7718.3
$ws.Range("J106").Value = 7718.3
$ws.Range("L106").Value = 7718.3
$ws.Range("N106").Value = -10242.3
$ws.Range("H122").Value = 6172.2104
$ws.Range("I122").Value = 4498.273
$ws.Range("K122").Value = 13494.819
$ws.Range("M122").Value = -11044.819
$ws.Range("H136").Value = 81879.5
$ws.Range("I136").Value = 169009.92
$ws.Range("J136").Value = 7196.2856
$ws.Range("K136").Value = 507029.76
$ws.Range("L136").Value = 21588.8568
$ws.Range("M136").Value = -504479.76
$ws.Range("N136").Value = -26688.8568
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 59999.6
$ws.Range("J4").Value = 59999.6
$ws.Range("L4").Value = 59999.6
$ws.Range("N4").Value = -60225.6
$ws.Range("H74").Value = 8309
$ws.Range("J74").Value = 8309
$ws.Range("L74").Value = 8309
$ws.Range("N74").Value = -10181
$ws.Range("H77").Value = 8309
$ws.Range("J77").Value = 8309
$ws.Range("L77").Value = 24927
$ws.Range("N77").Value = -34287
$ws.Range("H81").Value = 6671307.5
$ws.Range("I81").Value = 8775742
$ws.Range("K81").Value = 17551484
$ws.Range("M81").Value = -17550423
$ws.Range("H84").Value = 6671307.5
$ws.Range("I84").Value = 8775742
$ws.Range("K84").Value = 87757420
$ws.Range("M84").Value = -87752116
$ws.Range("H132").Value = 17437548
$ws.Range("I132").Value = 19615192
$ws.Range("K132").Value = 58845576
$ws.Range("M132").Value = -58843046
$ws.Range("H136").Value = 3905.0334
$ws.Range("I136").Value = 4255.5107
$ws.Range("J136").Value = 2637.923
$ws.Range("K136").Value = 12766.5321
$ws.Range("L136").Value = 7913.768999999999
$ws.Range("M136").Value = -10216.5321
$ws.Range("N136").Value = -13013.769
